# Quarterly dollar_cumulative update for darou/desobhan income statement.
# The 1402-02-28 estimated-publish-date quarter has been finalized/updated
# to an actual publish date of 1402-03-07, and the reported Q4 1401-12
# (column M) figures were refreshed with the finalized numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Publish-date header labels (row 9)
$ws.Range("I9").Value = "1402-03-07 (8)"
$ws.Range("M9").Value = "1402-03-07 (2)"

# Updated financial figures for column M (12 ماهه منتهی به 1401/12)
$ws.Range("M12").Value = -17794   # بهای تمام شده کالای فروش رفته
$ws.Range("M13").Value = 17218    # سود (زیان) ناخالص
$ws.Range("M14").Value = -1854    # هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("M15").Value = "-"      # هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی)
$ws.Range("M16").Value = 213      # خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("M17").Value = 15577    # سود (زیان) عملیاتی
$ws.Range("M19").Value = 259      # خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("M20").Value = 12970    # سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("M21").Value = -2128    # مالیات
$ws.Range("M22").Value = 10842    # سود (زیان) خالص عملیات در حال تداوم
$ws.Range("M24").Value = 10842    # سود (زیان) خالص
